$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug in probability calculator: certain edge-case transition
# probabilities were not being assigned correctly (missing the 0.875
# "stay" probability mass that should be added in boundary/edge cases).
$ws.Range("A28").Value = 0.95000000000000007
$ws.Range("A54").Value = 0.95000000000000007
$ws.Range("A288").Value = 0.95000000000000007
$ws.Range("A314").Value = 0.95000000000000007
$ws.Range("A522").Value = 0.95000000000000007
$ws.Range("A548").Value = 0.95000000000000007
$ws.Range("A574").Value = 0.95000000000000007
$ws.Range("A600").Value = 0.92500000000000004
$ws.Range("A626").Value = 0.95000000000000007
$ws.Range("A731").Value = 0.95000000000000007
$ws.Range("A757").Value = 0.97500000000000009
$ws.Range("A783").Value = 0.92500000000000004
$ws.Range("A861").Value = 0.92500000000000004
$ws.Range("A991").Value = 0.92500000000000004
$ws.Range("A1017").Value = 0.95000000000000007
$ws.Range("A1043").Value = 0.92500000000000004
$ws.Range("A1121").Value = 0.92500000000000004
$ws.Range("A1251").Value = 0.95000000000000007
$ws.Range("A1252").Value = 0.95000000000000007
$ws.Range("A1278").Value = 0.95000000000000007
$ws.Range("A1304").Value = 0.95000000000000007
$ws.Range("A1330").Value = 0.92500000000000004
$ws.Range("A1356").Value = 0.95000000000000007
$ws.Range("A1382").Value = 0.97500000000000009
$ws.Range("A1538").Value = 0.95000000000000007
$ws.Range("A1564").Value = 0.95000000000000007
$ws.Range("A1798").Value = 0.95000000000000007
$ws.Range("A1824").Value = 0.95000000000000007
$ws.Range("A1877").Value = 0.95000000000000007
$ws.Range("A2007").Value = 0.97500000000000009
$ws.Range("A2059").Value = 0.92500000000000004
$ws.Range("A2085").Value = 0.92500000000000004
$ws.Range("A2137").Value = 0.92500000000000004
$ws.Range("A2267").Value = 0.95000000000000007
$ws.Range("A2319").Value = 0.92500000000000004
$ws.Range("A2345").Value = 0.92500000000000004
$ws.Range("A2397").Value = 0.95000000000000007

# Update the active selection to reflect the used range, as Excel does
# after the user selects the full data range before saving.
$ws.Range("A1:A3126").Select()
